# ATAS.docx weekly update
# 1) "Emmily" / " Ferreira - Scrum Master" (split across two runs, with a
#    proofErr spell-check bookmark around "Emmily") becomes a single run
#    "Emmily Ferreira - Scrum Master" (11 occurrences across the weekly
#    tables).
# 2) Same fix for the rows where the text is just "Emmily" / " Ferreira "
#    (4 occurrences).
# 3) "alem" (flagged by the spell checker) corrected to "além".
# 4) "- Revisão " / "do banco de dados." merged into a single run.
# 5) "- " / "Definição dos responsáveis pelas novas telas." merged into a
#    single run.
# 6) The placeholder "=" run becomes "Emmily Ferreira".
# 7) The trailing " " run before "foi realizado a observação..." merges
#    into that run.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

$d.Content.Find.Execute(
    "Emmily Ferreira – Scrum Master", $false, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "Emmily Ferreira – Scrum Master", $wdReplaceAll)

$d.Content.Find.Execute(
    "Emmily Ferreira ", $false, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "Emmily Ferreira ", $wdReplaceAll)

$d.Content.Find.Execute(
    "alem", $true, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "além", $wdReplaceAll)

$d.Content.Find.Execute(
    "- Revisão do banco de dados.", $false, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "- Revisão do banco de dados.", $wdReplaceAll)

$d.Content.Find.Execute(
    "- Definição dos responsáveis pelas novas telas.", $false, $false,
    $false, $false, $false, $true, $wdFindContinue, $false,
    "- Definição dos responsáveis pelas novas telas.", $wdReplaceAll)

$d.Content.Find.Execute(
    "=", $false, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "Emmily Ferreira", $wdReplaceAll)

$d.Content.Find.Execute(
    "foi realizado a observação dos pontos a serem alterados no banco de dados, além disso a divisão de quem ficaria responsável por produzir cada tela nova (dashboard, cadastro de funcionário, armazém e o menu lateral que será utilizado nelas)",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "foi realizado a observação dos pontos a serem alterados no banco de dados, além disso a divisão de quem ficaria responsável por produzir cada tela nova (dashboard, cadastro de funcionário, armazém e o menu lateral que será utilizado nelas)",
    $wdReplaceAll)
